# khl/Injuries_Master_Clubs.xlsx - publish run (2025-12-01 03:0x UTC)
# Sheet "snapshot": two players recovered (removed) and one new injury added;
# Sheet "returned": log the two recovered players;
# Sheet "new_injured": log the newly injured player.

# Helper: write a value as plain text, even when it looks like a pure
# number ("20", "24799", ...) or a bare date ("2025-12-01") that Excel's
# smart-typing would otherwise coerce into a numeric/date cell. Formatting
# as Text first, then clearing the format again afterwards, keeps the
# cell's style index untouched (matches the source file, which has no
# explicit style on data cells).
function Set-TextValue($cell, $text) {
  $cell.NumberFormat = "@"
  $cell.Value = $text
  $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: snapshot
# ---------------------------------------------------------------
$snap = $wb.Worksheets.Item("snapshot")

# Players no longer injured (recovered / returned):
#   row 27 - ТОР Науменков Михаил
#   row 20 - СИБ Широков Сергей
# Delete the higher-numbered row first so the other row index stays valid.
$snap.Rows.Item(27).Delete()
$snap.Rows.Item(20).Delete()

# New injury to insert just before what is now row 26 (ЦСК Бучельников Дмитрий),
# i.e. right after СПР Вишневский Дмитрий (row 25) - keeps alphabetical team order.
$snap.Rows.Item(26).Insert()

$snap.Cells.Item(26,1).Value = "ТРК"
$snap.Cells.Item(26,2).Value = "Трактор"
$snap.Cells.Item(26,3).Value = "traktor"
$snap.Cells.Item(26,4).Value = "Мыльников Сергей И"
Set-TextValue $snap.Cells.Item(26,5) "20"
$snap.Cells.Item(26,6).Value = "вратарь"
Set-TextValue $snap.Cells.Item(26,7) "24799"
$snap.Cells.Item(26,8).Value = "1369_ТРК_мыльниковсергейи"
$snap.Cells.Item(26,9).Value = "injured_active"
$snap.Cells.Item(26,10).Value = "https://www.khl.ru/clubs/traktor/team/"
$snap.Cells.Item(26,11).Value = "2025-12-01T03:02:45.517776+00:00"

# Refresh scraped_at (column K) for every remaining/shifted data row to match
# the new scrape run.
$scrapedAt = @{
  2  = "2025-12-01T03:01:58.444613+00:00"
  3  = "2025-12-01T03:02:01.195064+00:00"
  4  = "2025-12-01T03:02:01.195094+00:00"
  5  = "2025-12-01T03:02:03.474917+00:00"
  6  = "2025-12-01T03:02:05.735149+00:00"
  7  = "2025-12-01T03:02:08.084028+00:00"
  8  = "2025-12-01T03:02:11.786477+00:00"
  9  = "2025-12-01T03:02:16.872729+00:00"
  10 = "2025-12-01T03:02:16.872757+00:00"
  11 = "2025-12-01T03:02:19.192996+00:00"
  12 = "2025-12-01T03:02:21.985029+00:00"
  13 = "2025-12-01T03:02:21.985060+00:00"
  14 = "2025-12-01T03:02:24.682199+00:00"
  15 = "2025-12-01T03:02:27.414139+00:00"
  16 = "2025-12-01T03:02:27.414169+00:00"
  17 = "2025-12-01T03:02:27.414186+00:00"
  18 = "2025-12-01T03:02:29.694643+00:00"
  19 = "2025-12-01T03:02:29.694674+00:00"
  20 = "2025-12-01T03:02:32.351899+00:00"
  21 = "2025-12-01T03:02:32.351937+00:00"
  22 = "2025-12-01T03:02:35.178916+00:00"
  23 = "2025-12-01T03:02:35.178948+00:00"
  24 = "2025-12-01T03:02:35.178977+00:00"
  25 = "2025-12-01T03:02:37.960343+00:00"
  27 = "2025-12-01T03:02:48.244707+00:00"
  28 = "2025-12-01T03:02:48.244743+00:00"
  29 = "2025-12-01T03:02:50.512767+00:00"
  30 = "2025-12-01T03:02:50.512794+00:00"
}
foreach ($r in $scrapedAt.Keys) {
  $snap.Cells.Item($r,11).Value = $scrapedAt[$r]
}

# ---------------------------------------------------------------
# Sheet 2: returned
# ---------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")

# Replace the previous run's two-return rows, drop the third (no longer returning).
$returned.Cells.Item(2,1).Value = "СИБ"
$returned.Cells.Item(2,2).Value = "Сибирь"
$returned.Cells.Item(2,3).Value = "Широков Сергей"
$returned.Cells.Item(2,4).Value = "1369_СИБ_широковсергей"
$returned.Cells.Item(2,5).Value = "RETURN"
$returned.Cells.Item(2,6).Value = "2025-12-01T11:02:51.019576+08:00"
Set-TextValue $returned.Cells.Item(2,7) "2025-12-01"

$returned.Cells.Item(3,1).Value = "ТОР"
$returned.Cells.Item(3,2).Value = "Торпедо"
$returned.Cells.Item(3,3).Value = "Науменков Михаил"
$returned.Cells.Item(3,4).Value = "1369_ТОР_науменковмихаил"
$returned.Cells.Item(3,5).Value = "RETURN"
$returned.Cells.Item(3,6).Value = "2025-12-01T11:02:51.019576+08:00"
Set-TextValue $returned.Cells.Item(3,7) "2025-12-01"

$returned.Rows.Item(4).Delete()

# ---------------------------------------------------------------
# Sheet 3: new_injured
# ---------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")

$newInjured.Cells.Item(2,1).Value = "ТРК"
$newInjured.Cells.Item(2,2).Value = "Трактор"
$newInjured.Cells.Item(2,3).Value = "Мыльников Сергей И"
$newInjured.Cells.Item(2,4).Value = "1369_ТРК_мыльниковсергейи"
$newInjured.Cells.Item(2,5).Value = "INJURED_NEW"
$newInjured.Cells.Item(2,6).Value = "2025-12-01T11:02:51.019576+08:00"
Set-TextValue $newInjured.Cells.Item(2,7) "2025-12-01"

Write-Output "edit applied"
